$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Role (Principal/Sub)" column header and its values to the
# new "Role (PI/Sub I)" wording (Principal -> PI, sub -> sub I).
$ws.Range("A1").Value = "Role (PI/Sub I)"
$ws.Range("A2").Value = "PI"
$ws.Range("A3").Value = "sub I"
$ws.Range("A4").Value = "PI"
$ws.Range("A5").Value = "sub I"
$ws.Range("A6").Value = "sub I"
$ws.Range("A7").Value = "sub I"

# Move the active selection to A8, matching the saved view state.
$ws.Range("A8").Select()
